$wb = $excel.ActiveWorkbook

$wsTDIL = $wb.Worksheets.Item("TDIL-EN")
$wsSPDS = $wb.Worksheets.Item("SPDS-JP")

# --- TDIL-EN: insert a new row 2 for "Magical Something" ---------------
# Shift existing data rows (2..10) down to (3..11) by copying values,
# bottom-up so we don't clobber data before it's copied.
for ($r = 10; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $wsTDIL.Range("A$dst").Value = $wsTDIL.Range("A$src").Value2
    $wsTDIL.Range("B$dst").Value = $wsTDIL.Range("B$src").Value2
    $wsTDIL.Range("C$dst").Value = $wsTDIL.Range("C$src").Value2
    $wsTDIL.Range("D$dst").Value = $wsTDIL.Range("D$src").Value2
    $wsTDIL.Range("E$dst").Value = $wsTDIL.Range("E$src").Value2
}

# New row 2 content
$wsTDIL.Range("A2").Value = "Magical Something"
$wsTDIL.Range("B2").Value = 100909000

# --- SPDS-JP: add name for row 22 ---------------------------------------
$wsSPDS.Range("A22").Value = "Abyss Actor - Pretty Heroine"

# --- selections / active sheet -----------------------------------------
$wsTDIL.Range("F5").Select()
$wsSPDS.Range("F20").Select()
$wsSPDS.Activate()
